# Auto-generated edit script applying the Spriggan_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 16265.066
$ws.Range("J69").Value = 17229.924
$ws.Range("L69").Value = 51689.772
$ws.Range("N69").Value = -53437.772
$ws.Range("H70").Value = 11150
$ws.Range("I70").Value = 4214.2856
$ws.Range("K70").Value = 12642.8568
$ws.Range("M70").Value = -12372.8568
$ws.Range("H72").Value = 16265.066
$ws.Range("J72").Value = 17229.924
$ws.Range("L72").Value = 155069.316
$ws.Range("N72").Value = -163805.316
$ws.Range("H73").Value = 11150
$ws.Range("I73").Value = 4214.2856
$ws.Range("K73").Value = 12642.8568
$ws.Range("M73").Value = -11706.8568
$ws.Range("H132").Value = 2521.5925
$ws.Range("I132").Value = 2576.2693
$ws.Range("K132").Value = 7728.8079
$ws.Range("M132").Value = -5198.8079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 430265.06
$ws.Range("I2").Value = 981735.4
$ws.Range("K2").Value = 981735.4
$ws.Range("M2").Value = -981622.4
$ws.Range("H32").Value = 6997.297
$ws.Range("I32").Value = 3752.3264
$ws.Range("J32").Value = 17597.533
$ws.Range("K32").Value = 3752.3264
$ws.Range("L32").Value = 17597.533
$ws.Range("M32").Value = -3465.3264
$ws.Range("N32").Value = -18171.533
$ws.Range("H45").Value = 981.6
$ws.Range("I45").Value = 856.53845
$ws.Range("K45").Value = 856.53845
$ws.Range("M45").Value = -479.53845
$ws.Range("H116").Value = 430265.06
$ws.Range("I116").Value = 981735.4
$ws.Range("K116").Value = 981735.4
$ws.Range("M116").Value = -979441.4
$ws.Range("H132").Value = 1700754.2
$ws.Range("I132").Value = 2005621.1
$ws.Range("J132").Value = 7049.5557
$ws.Range("K132").Value = 6016863.300000001
$ws.Range("L132").Value = 21148.6671
$ws.Range("M132").Value = -6014333.300000001
$ws.Range("N132").Value = -26208.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 430265.06
$ws.Range("I3").Value = 981735.4
$ws.Range("K3").Value = 981735.4
$ws.Range("M3").Value = -981621.4
$ws.Range("H18").Value = 100000
$ws.Range("J18").Value = 100000
$ws.Range("L18").Value = 100000
$ws.Range("N18").Value = -101058
$ws.Range("H134").Value = 12502002
$ws.Range("I134").Value = 12502002
$ws.Range("K134").Value = 37506006
$ws.Range("M134").Value = -37503471

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 794.5
$ws.Range("I29").Value = 794.5
$ws.Range("K29").Value = 794.5
$ws.Range("M29").Value = -501.5
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H107").Value = 635667.9
$ws.Range("I107").Value = 679574.4399999999
$ws.Range("J107").Value = 401499.66
$ws.Range("K107").Value = 679574.4399999999
$ws.Range("L107").Value = 401499.66
$ws.Range("M107").Value = -677654.4399999999
$ws.Range("N107").Value = -405339.66
$ws.Range("H134").Value = 20836166
$ws.Range("I134").Value = 25002108
$ws.Range("J134").Value = 6457
$ws.Range("K134").Value = 75006324
$ws.Range("L134").Value = 19371
$ws.Range("M134").Value = -75003789
$ws.Range("N134").Value = -24441

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 64549.312
$ws.Range("I5").Value = 91825
$ws.Range("J5").Value = 4542.8
$ws.Range("K5").Value = 275475
$ws.Range("L5").Value = 13628.4
$ws.Range("M5").Value = -275363
$ws.Range("N5").Value = -13852.4
$ws.Range("H135").Value = 64549.312
$ws.Range("I135").Value = 91825
$ws.Range("J135").Value = 4542.8
$ws.Range("K135").Value = 826425
$ws.Range("L135").Value = 40885.2
$ws.Range("M135").Value = -823890
$ws.Range("N135").Value = -45955.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 500
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -726
$ws.Range("H70").Value = 9836.75
$ws.Range("I70").Value = 12799.2
$ws.Range("K70").Value = 12799.2
$ws.Range("M70").Value = -12529.2
$ws.Range("H73").Value = 9836.75
$ws.Range("I73").Value = 12799.2
$ws.Range("K73").Value = 12799.2
$ws.Range("M73").Value = -11863.2
$ws.Range("H122").Value = 68516.664
$ws.Range("I122").Value = 93792.38
$ws.Range("K122").Value = 281377.14
$ws.Range("M122").Value = -278927.14

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1703.8
$ws.Range("I22").Value = 2084.875
$ws.Range("J22").Value = 1268.2858
$ws.Range("K22").Value = 2084.875
$ws.Range("L22").Value = 1268.2858
$ws.Range("M22").Value = -1789.875
$ws.Range("N22").Value = -1858.2858
$ws.Range("H27").Value = 1703.8
$ws.Range("I27").Value = 2084.875
$ws.Range("J27").Value = 1268.2858
$ws.Range("K27").Value = 2084.875
$ws.Range("L27").Value = 1268.2858
$ws.Range("M27").Value = -1977.875
$ws.Range("N27").Value = -1482.2858
$ws.Range("H40").Value = 13357.4
$ws.Range("I40").Value = 11104.692
$ws.Range("J40").Value = 28000
$ws.Range("K40").Value = 11104.692
$ws.Range("L40").Value = 28000
$ws.Range("M40").Value = -10968.692
$ws.Range("N40").Value = -28272
$ws.Range("H122").Value = 2835.4211
$ws.Range("I122").Value = 2138.4546
$ws.Range("K122").Value = 6415.3638
$ws.Range("M122").Value = -3965.3638
$ws.Range("H136").Value = 2999
$ws.Range("J136").Value = 2999
$ws.Range("L136").Value = 8997
$ws.Range("N136").Value = -14097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 23760.77
$ws.Range("J54").Value = 23749.25
$ws.Range("L54").Value = 23749.25
$ws.Range("N54").Value = -24789.25
$ws.Range("H81").Value = 986.2857
$ws.Range("I81").Value = 935
$ws.Range("J81").Value = 1024.75
$ws.Range("K81").Value = 1870
$ws.Range("L81").Value = 2049.5
$ws.Range("M81").Value = -809
$ws.Range("N81").Value = -4171.5
$ws.Range("H84").Value = 986.2857
$ws.Range("I84").Value = 935
$ws.Range("J84").Value = 1024.75
$ws.Range("K84").Value = 9350
$ws.Range("L84").Value = 10247.5
$ws.Range("M84").Value = -4046
$ws.Range("N84").Value = -20855.5
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 5975.6665
$ws.Range("I122").Value = 5975.6665
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 17926.9995
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15476.9995
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 38464732
$ws.Range("I136").Value = 55557524
$ws.Range("K136").Value = 166672572
$ws.Range("M136").Value = -166670022
